# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (interest-count) column F across all four
# sheets to the freshly scraped numbers.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1770
$ws.Range("F5").Value = 537
$ws.Range("F6").Value = 2146
$ws.Range("F8").Value = 2093
$ws.Range("F11").Value = 2413
$ws.Range("F12").Value = 663
$ws.Range("F14").Value = 3960
$ws.Range("F16").Value = 370
$ws.Range("F17").Value = 3058
$ws.Range("F18").Value = 830
$ws.Range("F19").Value = 145
$ws.Range("F21").Value = 140
$ws.Range("F22").Value = 2066
$ws.Range("F23").Value = 1185
$ws.Range("F24").Value = 1919
$ws.Range("F25").Value = 394
$ws.Range("F26").Value = 209
$ws.Range("F27").Value = 18
$ws.Range("F28").Value = 8531
$ws.Range("F29").Value = 5742
$ws.Range("F30").Value = 355
$ws.Range("F31").Value = 178
$ws.Range("F32").Value = 755
$ws.Range("F33").Value = 772
$ws.Range("F34").Value = 3478
$ws.Range("F37").Value = 396
$ws.Range("F38").Value = 40
$ws.Range("F39").Value = 194
$ws.Range("F41").Value = 4646
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 863
$ws.Range("F44").Value = 80

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 6
$ws.Range("F13").Value = 20
$ws.Range("F18").Value = 421

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8323
$ws.Range("F4").Value = 1296

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1296
$ws.Range("F6").Value = 1770
$ws.Range("F8").Value = 537
$ws.Range("F10").Value = 2093
$ws.Range("F15").Value = 3960
$ws.Range("F16").Value = 370
$ws.Range("F17").Value = 3058
$ws.Range("F18").Value = 830
$ws.Range("F19").Value = 145
$ws.Range("F21").Value = 2066
$ws.Range("F25").Value = 20
$ws.Range("F27").Value = 1919
$ws.Range("F29").Value = 209
$ws.Range("F30").Value = 18
$ws.Range("F31").Value = 8531
$ws.Range("F32").Value = 5742
$ws.Range("F34").Value = 355
$ws.Range("F35").Value = 178
$ws.Range("F36").Value = 755
$ws.Range("F37").Value = 772
$ws.Range("F40").Value = 396
$ws.Range("F41").Value = 194
$ws.Range("F43").Value = 4646
$ws.Range("F44").Value = 863
